# Regenerate Report for Handback
# Replace old source-file UUID-based names/timestamps with the new ones,
# on all three sheets (Overview, zh-cn, de-de), and keep the hyperlink
# "display" text in sync with the new file names.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# New identifiers / timestamps
# ---------------------------------------------------------------------
$oldFile1 = "1f73811e-a810-4e5f-a9db-ff849dfa4bfa"
$newFile1 = "daad096a-e73d-4d74-9941-51217f3a201d"
$oldFile2 = "ffdd4f67-b553-47d7-a0ae-7a81ec919525"
$newFile2 = "ffff7d060d44-1892-4b94-bc38-417507093462"

$newFile1Md = "$newFile1.md"
$newFile2Md = "$newFile2.md"
$newFile1MdPath = "e2e\$newFile1.md"
$newFile2MdPath = "e2e\$newFile2.md"

$newXliffZhCn = "$newFile1.4be0fcbe55b7edc038c0970a1b86a9643750baab.zh-cn.xlf"
$newXliffDeDe = "$newFile1.4be0fcbe55b7edc038c0970a1b86a9643750baab.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newFile1Md
$wsOverview.Range("B2").Value = $newFile1MdPath
$wsOverview.Range("G2").Value = "2016-08-21 13:06:17"

$wsOverview.Range("A3").Value = $newFile2Md
$wsOverview.Range("B3").Value = $newFile2MdPath
$wsOverview.Range("G3").Value = "2016-08-21 13:06:17"

# Update hyperlink display text for B2 / B3 while preserving the original
# relationship targets (Address) and ordering (rId2, rId3).
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B3").Hyperlinks.Delete()

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1MdPath) | Out-Null

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2MdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn.Range("A2").Value = $newFile1Md
$wsZhCn.Range("G2").Value = $newXliffZhCn
$wsZhCn.Range("H2").Value = "2016-08-21 13:06:13"
$wsZhCn.Range("I2").Value = $newFile1Md
$wsZhCn.Range("J2").Value = $newXliffZhCn
$wsZhCn.Range("K2").Value = "2016-08-21 13:06:29"

$wsZhCn.Range("A3").Value = $newFile2Md
$wsZhCn.Range("G3").Value = $newXliffZhCn
$wsZhCn.Range("H3").Value = "2016-08-21 13:06:13"
$wsZhCn.Range("I3").Value = $newFile2Md
$wsZhCn.Range("J3").Value = $newXliffZhCn
$wsZhCn.Range("K3").Value = "2016-08-21 13:06:29"

# Update hyperlink display text for A2 / I2 / A3 / I3, preserving the
# original relationship targets (Address) and ordering (rId2..rId5).
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Range("I3").Hyperlinks.Delete()

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/18c72ac7779159390235baa950209fe859c93f87/e2e/$oldFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/18c72ac7779159390235baa950209fe859c93f87/e2e/$oldFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md) | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe.Range("A2").Value = $newFile1Md
$wsDeDe.Range("G2").Value = $newXliffDeDe
$wsDeDe.Range("H2").Value = "2016-08-21 13:06:17"
$wsDeDe.Range("I2").Value = $newFile1Md
$wsDeDe.Range("J2").Value = $newXliffDeDe
$wsDeDe.Range("K2").Value = "2016-08-21 13:06:35"

$wsDeDe.Range("A3").Value = $newFile2Md
$wsDeDe.Range("G3").Value = $newXliffDeDe
$wsDeDe.Range("H3").Value = "2016-08-21 13:06:17"
$wsDeDe.Range("I3").Value = $newFile2Md
$wsDeDe.Range("J3").Value = $newXliffDeDe
$wsDeDe.Range("K3").Value = "2016-08-21 13:06:35"

# Update hyperlink display text for A2 / I2 / A3 / I3, preserving the
# original relationship targets (Address) and ordering (rId2..rId5).
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Range("I3").Hyperlinks.Delete()

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/15b70300accc379a712c0ba663f46cb78f41ebc0/e2e/$oldFile1.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile1Md) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af2ef8531f08e5c095f3672bf4ddc2c8402bbe7a/e2e/$oldFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/15b70300accc379a712c0ba663f46cb78f41ebc0/e2e/$oldFile2.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    $newFile2Md) | Out-Null

Write-Host "Handback status report regenerated."
